$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns P1, Q1 (continuing the 0..n sequence), copying the
# formatting of the existing header cell O1 (bold, centered, bordered)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2-25), update columns I-O with the new values and add
# the new P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I
    $ws.Cells.Item($r, 10).Value = 2   # J
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 12).Value = 2   # L
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 14).Value = 2   # N
    $ws.Cells.Item($r, 15).Value = 1   # O
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
